# Edit LOM3228.docx: split the "Programa" paragraphs (PT and EN) with a
# manual line break (<w:br/>) right before "Criogenia." / "cryogenics."
# respectively, per the commit diff.

$d = $word.ActiveDocument

$orig1 = 'Teoria dos gases rarefeitos. Escoamento de gases. Bombas de vácuo. Descrição quantitativa do bombeamento de sistemas de vácuo. Medidores de pressão. Acessórios: armadilhas, anteparos, válvulas, etc. Adsorção, dessorção e evaporação de moléculas em vácuo. Detecção de vazamento. Vedação. Soldagem. Limpeza. Criogenia. Propriedades de gases e líquidos criogênicos. Métodos para obtenção de baixa temperatura. Liquefação de gases. Medição de temperatura. Componentes criogênicos. Cálculo de transferência de calor em criostatos e dewars.'
$orig2 = 'Theory of rarefied gases. Gas flow. Vacuum pumps. Quantitative description of the pumping of vacuum systems. Pressure gauges. Accessories: traps, shields, valves, etc. Adsorption, desorption and evaporation of molecules in vacuum. Leak detection .Sealing.Welding.Cleaning.cryogenics. Properties of cryogenic gases and liquids. Methods for obtaining low temperature. Liquefaction of gases. Temperature measurement. Cryogenic components. Calculation of heat transfer in cryostats and dewars.'

$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Teoria dos gases rarefeitos. Escoamento de gases. Bombas de vácuo. Descrição quantitativa do bombeamento de sistemas de vácuo. Medidores de pressão. Acessórios: armadilhas, anteparos, válvulas, etc. Adsorção, dessorção e evaporação de moléculas em vácuo. Detecção de vazamento. Vedação. Soldagem. Limpeza. </w:t><w:br/><w:t>Criogenia. Propriedades de gases e líquidos criogênicos. Métodos para obtenção de baixa temperatura. Liquefação de gases. Medição de temperatura. Componentes criogênicos. Cálculo de transferência de calor em criostatos e dewars.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Theory of rarefied gases. Gas flow. Vacuum pumps. Quantitative description of the pumping of vacuum systems. Pressure gauges. Accessories: traps, shields, valves, etc. Adsorption, desorption and evaporation of molecules in vacuum. Leak detection .Sealing.Welding.Cleaning.</w:t><w:br/><w:t>cryogenics. Properties of cryogenic gases and liquids. Methods for obtaining low temperature. Liquefaction of gases. Temperature measurement. Cryogenic components. Calculation of heat transfer in cryostats and dewars.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

function Find-ParagraphIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        # Trim the trailing paragraph-mark / cell-mark characters Word
        # appends to Range.Text so we compare on content only.
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

$idx1 = Find-ParagraphIndexByText $d $orig1
if ($idx1 -lt 0) {
    throw "Could not locate the Portuguese 'Programa' paragraph"
}
$r1 = $d.Paragraphs.Item($idx1).Range
$null = $r1.InsertXML($xml1)

$idx2 = Find-ParagraphIndexByText $d $orig2
if ($idx2 -lt 0) {
    throw "Could not locate the English 'Programa' paragraph"
}
$r2 = $d.Paragraphs.Item($idx2).Range
$null = $r2.InsertXML($xml2)

Write-Output "OK: updated paragraphs $idx1 and $idx2"
